# Updated legacy GSC export data:
# The export window advanced by one day: the oldest day (2025-10-09) drops
# off the top of the table, every remaining day shifts up one row, and the
# three most-recent days (which Search Console hasn't finished indexing/
# reporting on yet) now have blank "No video indexed" / "Video indexed"
# counts instead of numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Drop the oldest row (2025-10-09); everything below shifts up one row.
$ws.Range("A2").EntireRow.Delete()

# The newest three days (now rows 2-4) have no data yet.
$ws.Range("B2:C4").Value = ""
